$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mis-keyed record: customer id 1004 -> 1007 for the
# "Test Wrong ID" row (currently row 8, column A).
$ws.Range("A8").Value = 1007

# Re-sort the data range A2:D10 by column A (customr_id) ascending,
# just like Data > Sort would do after correcting the id above.
$ws.Sort.SortFields.Clear()
$sortRange = $ws.Range("A2:D10")
$keyRange = $ws.Range("A2")
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Format column A (the id column) as a plain integer number format and
# select the whole column, as the author apparently did after sorting.
$ws.Columns.Item(1).NumberFormat = "0"
$ws.Columns.Item(1).Select()

# Make sure the page is set to portrait orientation.
$ws.PageSetup.Orientation = 1
